$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad") for all existing data rows (2-325) moves from
#    45202 to 45203 (one day later).
$ws.Range("C2:C325").Value = 45203

# 2) Two new cases were appended at the bottom of the sheet: rows 326 & 327.
#    They follow the same layout as the other "Sveaskog" rows (e.g. row 325):
#    columns A-Q populated, R present but empty (wrap-text style), no
#    hyperlink formulas in S:Y.

# --- Row 326 : A 47300-2023 ---
$ws.Cells.Item(326, 1).Value = "A 47300-2023"
$ws.Cells.Item(326, 2).Value = 45202
$ws.Cells.Item(326, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(326, 3).Value = 45203
$ws.Cells.Item(326, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(326, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(326, 5).Value = "ANEBY"
$ws.Cells.Item(326, 6).Value = "Sveaskog"
$ws.Cells.Item(326, 7).Value = 1.8
$ws.Cells.Item(326, 8).Value = 0
$ws.Cells.Item(326, 9).Value = 0
$ws.Cells.Item(326, 10).Value = 0
$ws.Cells.Item(326, 11).Value = 0
$ws.Cells.Item(326, 12).Value = 0
$ws.Cells.Item(326, 13).Value = 0
$ws.Cells.Item(326, 14).Value = 0
$ws.Cells.Item(326, 15).Value = 0
$ws.Cells.Item(326, 16).Value = 0
$ws.Cells.Item(326, 17).Value = 0
$ws.Cells.Item(326, 18).Value = ""
$ws.Cells.Item(326, 18).WrapText = $true
$ws.Rows.Item(326).RowHeight = 15

# --- Row 327 : A 47296-2023 ---
$ws.Cells.Item(327, 1).Value = "A 47296-2023"
$ws.Cells.Item(327, 2).Value = 45202
$ws.Cells.Item(327, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(327, 3).Value = 45203
$ws.Cells.Item(327, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(327, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(327, 5).Value = "ANEBY"
$ws.Cells.Item(327, 6).Value = "Sveaskog"
$ws.Cells.Item(327, 7).Value = 0.6
$ws.Cells.Item(327, 8).Value = 0
$ws.Cells.Item(327, 9).Value = 0
$ws.Cells.Item(327, 10).Value = 0
$ws.Cells.Item(327, 11).Value = 0
$ws.Cells.Item(327, 12).Value = 0
$ws.Cells.Item(327, 13).Value = 0
$ws.Cells.Item(327, 14).Value = 0
$ws.Cells.Item(327, 15).Value = 0
$ws.Cells.Item(327, 16).Value = 0
$ws.Cells.Item(327, 17).Value = 0
$ws.Cells.Item(327, 18).Value = ""
$ws.Cells.Item(327, 18).WrapText = $true

# 3) Row 325 gains an explicit row height (matches the rest of the sheet).
$ws.Rows.Item(325).RowHeight = 15
